$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 294, pushing the existing row 294 (and all
# rows below it) down by one. This mirrors the existing data's row layout
# (D column carries a date-formatted style) since Insert() extends the
# style of the row above.
$ws.Rows(294).Insert()

# Populate the newly inserted row 294 with the new weekly price record.
$ws.Cells.Item(294, 1).Value2  = 5
$ws.Cells.Item(294, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(294, 3).Value2  = "Maule"
$ws.Cells.Item(294, 4).Value2  = 44943
$ws.Cells.Item(294, 5).Value2  = 7
$ws.Cells.Item(294, 6).Value2  = 100112009
$ws.Cells.Item(294, 7).Value2  = "Acelga"
$ws.Cells.Item(294, 8).Value2  = "Sin especificar"
$ws.Cells.Item(294, 9).Value2  = "Primera"
$ws.Cells.Item(294, 10).Value2 = 800
$ws.Cells.Item(294, 11).Value2 = 3000
$ws.Cells.Item(294, 12).Value2 = 3000
$ws.Cells.Item(294, 13).Value2 = 3000
$ws.Cells.Item(294, 14).Value2 = "$/docena de atados (4 kilos)"
$ws.Cells.Item(294, 15).Value2 = "Región del Maule"
$ws.Cells.Item(294, 16).Value2 = 750
$ws.Cells.Item(294, 17).Value2 = 4
$ws.Cells.Item(294, 18).Value2 = "Hortaliza"
